$d = $word.ActiveDocument

# Locate the first paragraph's run containing the original sentence and
# grab a Range right after it so we can append the new, differently
# colored run without disturbing the existing text/formatting.
$para = $d.Paragraphs.First
$r = $para.Range

# Trim the trailing paragraph mark off the paragraph range so our
# insertions land right after the existing text, before the pilcrow.
$r.End = $r.End - 1

# Append two trailing spaces to the existing sentence (kept in the
# original run's formatting).
$r.InsertAfter("  ")

# Now insert the new colored run right after that, in its own Range so
# we can give it distinct (red) formatting without affecting the rest.
$newRange = $d.Range($r.End, $r.End)
$newRange.InsertAfter([char]0x0028 + "This is a change " + [char]0x2013 + " Version for branch alternate" + [char]0x0029)
$newRange.Font.Color = 192
